$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.060.08'
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = '1.910.32'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4600'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -1.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07753'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9659'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '1.929.65'
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.013'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.722'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.07075'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '84.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009584'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.003'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.62%  '
$ws.Range("D22").Value = '29.060.77'
$ws.Range("E22").Value = '  -0.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.440'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.92%  '
$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.154.26'
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.095'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.17%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.667'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.36%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.59'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("B31").Value = 'LidoDAOToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.821'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.90%  '
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09302'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.19%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8551'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.53%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.104'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.259'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.43%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.079'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.160'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05680'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02050'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.509'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.04%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5563'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.13%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1760'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000002911'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.41%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.222'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.02%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.712'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.30%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5209'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.21%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06810'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.77%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.058'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.50%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.789'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.63'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.99%  '
